$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update price (D) / volume (E) figures. Several "Price" strings look
# numeric (e.g. "0.9998", "1.000", "0.000007481") but the source file
# stores them as plain inline-string text, so mark those cells as Text
# before assigning, otherwise Excel auto-converts them to numbers and
# mangles the formatting (trailing zeros, scientific notation, etc).

$ws.Range("D2").Value = "30.483.18"
$ws.Range("E2").Value = "  -1.27%  "
$ws.Range("D3").Value = "1.911.83"
$ws.Range("E3").Value = "  -1.43%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9998"
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "239.87"
$ws.Range("E5").Value = "  -1.59%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9992"
$ws.Range("E6").Value = "  -0.12%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4778"
$ws.Range("E7").Value = "  -2.76%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2844"
$ws.Range("E8").Value = "  -3.42%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06699"
$ws.Range("E9").Value = "  -2.80%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.49"
$ws.Range("E10").Value = "  +1.23%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "103.05"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07769"
$ws.Range("D13").Value = "1.911.19"
$ws.Range("E13").Value = "  -1.25%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.197"
$ws.Range("E14").Value = "  -3.19%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6709"
$ws.Range("E15").Value = "  -4.71%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "275.14"
$ws.Range("E16").Value = "  -0.08%  "
$ws.Range("D17").Value = "30.460.90"
$ws.Range("E17").Value = "  -1.33%  "
$ws.Range("E18").Value = "  -0.19%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007481"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.63"
$ws.Range("E20").Value = "  -3.70%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.391"
$ws.Range("E21").Value = "  -3.36%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.4635"
$ws.Range("E22").Value = "  -7.51%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.000"
$ws.Range("E23").Value = "  +0.08%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.305"
$ws.Range("E24").Value = "  -3.69%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.354"
$ws.Range("E25").Value = "  -5.14%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "167.93"
$ws.Range("E26").Value = "  +1.04%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "19.22"
$ws.Range("E27").Value = "  -2.13%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.083"
$ws.Range("E28").Value = "  -3.44%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.381"
$ws.Range("E29").Value = "  -0.79%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.09985"
$ws.Range("E30").Value = "  -4.45%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.577"
$ws.Range("E31").Value = "  +0.23%  "
$ws.Range("E32").Value = "  -3.17%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.263"
$ws.Range("E33").Value = "  -2.81%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.04735"
$ws.Range("E34").Value = "  -3.30%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7264"
$ws.Range("E35").Value = "  -4.28%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.114"
$ws.Range("E36").Value = "  -3.39%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.716"
$ws.Range("E37").Value = "  -0.76%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01905"
$ws.Range("E38").Value = "  -5.18%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.620"
$ws.Range("E39").Value = "  -1.49%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.333"
$ws.Range("E40").Value = "  -2.88%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "74.00"
$ws.Range("E41").Value = "  -5.29%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.8612"
$ws.Range("E44").Value = "  -5.71%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.4262"
$ws.Range("E45").Value = "  -4.10%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.9985"
$ws.Range("E46").Value = "  -0.03%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "7.413"
$ws.Range("E47").Value = "  -3.45%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "952.97"
$ws.Range("E48").Value = "  -4.46%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.1204"
$ws.Range("E49").Value = "  -3.49%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "34.68"
$ws.Range("E50").Value = "  -4.02%  "

# Rows 42/43 swap Quant <-> RenderToken (name, link, price, volume);
# row 51 replaces Cronos with EnergySwap.
$ws.Range("B42").Value = "RenderToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.962"
$ws.Range("E42").Value = "  -6.43%  "
$ws.Range("B43").Value = "Quant"
$ws.Range("C43").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "106.81"
$ws.Range("E43").Value = "  -0.80%  "
$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "8.794"
$ws.Range("E51").Value = "  -4.51%  "
